$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need NumberFormat forced to
# Text ("@") first, so Excel keeps the exact original string (with trailing
# zeros etc.) instead of silently converting them to floating point numbers.
$ws.Range('D2').Value = '63.713.53'
$ws.Range('D3').Value = '3.085.62'
$ws.Range('E3').Value = '  -2.26%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('E5').Value = '  -0.37%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '153.20'
$ws.Range('E6').Value = '  +4.00%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('E8').Value = '  +0.37%  '
$ws.Range('D9').Value = '3.083.24'
$ws.Range('E9').Value = '  -2.20%  '
$ws.Range('E10').Value = '  -2.49%  '
$ws.Range('E11').Value = '  -0.73%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.459'
$ws.Range('E12').Value = '  -0.83%  '
$ws.Range('B13').Value = 'Avalanche'
$ws.Range('C13').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '37.77'
$ws.Range('E13').Value = '  +1.47%  '
$ws.Range('B14').Value = 'ShibaInu'
$ws.Range('C14').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000243'
$ws.Range('E14').Value = '  -2.31%  '
$ws.Range('D15').Value = '3.597.25'
$ws.Range('E15').Value = '  -2.24%  '
$ws.Range('E16').Value = '  -2.12%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '7.16'
$ws.Range('E17').Value = '  -0.90%  '
$ws.Range('D18').Value = '63.636.90'
$ws.Range('E18').Value = '  -0.80%  '
$ws.Range('D19').Value = '3.081.34'
$ws.Range('E19').Value = '  -2.27%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '470.71'
$ws.Range('E20').Value = '  +0.52%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.66'
$ws.Range('E21').Value = '  +1.39%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.728'
$ws.Range('E22').Value = '  -1.20%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.56'
$ws.Range('E23').Value = '  +0.71%  '
$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '13.26'
$ws.Range('E24').Value = '  +1.50%  '
$ws.Range('B25').Value = 'Fetch.AI'
$ws.Range('C25').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.38'
$ws.Range('E25').Value = '  +1.23%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '81.24'
$ws.Range('E26').Value = '  -0.22%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.89'
$ws.Range('E28').Value = '  +1.25%  '
$ws.Range('E29').Value = '  -1.36%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.32'
$ws.Range('E30').Value = '  -0.77%  '
$ws.Range('E31').Value = '  +0.20%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.21'
$ws.Range('E32').Value = '  -0.74%  '
$ws.Range('E33').Value = '  +4.32%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '27.38'
$ws.Range('E34').Value = '  -0.54%  '
$ws.Range('D35').Value = '0.0₃0852'
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('E36').Value = '  -1.07%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.42'
$ws.Range('E37').Value = '  +4.17%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '6.13'
$ws.Range('E38').Value = '  -0.05%  '
$ws.Range('E39').Value = '  -4.23%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '9.36'
$ws.Range('E40').Value = '  +2.43%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '50.77'
$ws.Range('E41').Value = '  -2.35%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '447.04'
$ws.Range('E42').Value = '  -2.25%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.289'
$ws.Range('E43').Value = '  -1.76%  '
$ws.Range('E44').Value = '  -2.13%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '40.38'
$ws.Range('E45').Value = '  -0.20%  '
$ws.Range('D46').Value = '2.834.17'
$ws.Range('E46').Value = '  -3.36%  '
$ws.Range('E47').Value = '  -0.19%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '129.48'
$ws.Range('E48').Value = '  +1.20%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '25.56'
$ws.Range('E49').Value = '  +4.21%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.27'
$ws.Range('E51').Value = '  +1.13%  '
